# LoginCredential.xlsx - add a "Sr no" column in front of the teststeps data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("teststeps")

# Insert a new first column; emailid/password data shifts from A/B to B/C.
$ws.Columns.Item(1).Insert()

# Give the new column's header/value cells the same plain formatting used
# elsewhere on the row (copy an existing cell rather than re-deriving the
# format, so no extra style gets appended to styles.xml).
$ws.Range("B1").Copy()
$ws.Range("A1:A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header + serial number value
$ws.Range("A1").Value = "Sr no"
$ws.Range("A2").Value = 1

# The hyperlinks collection keeps pointing at its original cells after a
# column insert, so re-target it: stash the (still correct) hyperlink
# formatting from B2, delete + recreate the links against their new
# (shifted) cells, then restore the original formatting over the Add
# call's own formatting.
$fmtHolder = $ws.Range("Z1000")
$ws.Range("B2").Copy()
$fmtHolder.PasteSpecial(-4122)

$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:thoranere@rknec.edu")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Rajani@1992")

$fmtHolder.Copy()
$ws.Range("B2:C2").PasteSpecial(-4122)
$fmtHolder.Clear()
$excel.CutCopyMode = 0

# Column widths for the new layout
$ws.Columns.Item(1).ColumnWidth = 8.85546875
$ws.Columns.Item(2).ColumnWidth = 20.85546875
$ws.Columns.Item(3).ColumnWidth = 19.7109375

# Selection as recorded in the saved view state
$ws.Range("C7").Select()

$wb.Save()
